$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BME_DI_SPH")
$ws.Activate()

# --- Column B (actual_partial_oh): was C - D, now C - D - E ---
$ws.Range("B2").Formula = "=C2-D2-E2"
$ws.Range("B4:B8").Formula = "=C4-D4-E4"
$ws.Range("B3").Formula = "=C3-D3-E3"

# --- Column G (budgeted_partial_oh): was H - I, now H - I - J ---
$ws.Range("G2").Formula = "=H2-I2-J2"
$ws.Range("G3:G8").Formula = "=H3-I3-J3"

# --- Window / view state ---
$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
$win.Width = 29040
$win.Height = 15840

# Move the selection to G12 (also clears the stale topLeftCell="B1" hint).
[void]$ws.Range("G12").Select()
